$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update Status column (B) from "Failed" to "Passed" for the rows that
#        are now passing. Use PasteSpecial(formats) from B2 (an existing
#        "Passed" cell) so the already-defined bold+green style (s=4) is
#        reused instead of minting a near-duplicate style.
$ws.Range("B2").Copy()
$passedRows = @(3,5,6,7,8,9,10,17)
foreach ($r in $passedRows) {
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("B$r").Value = "Passed"
}

# --- 2. Update Run No. column (D)
$ws.Range("D10").Value = 3
$runNoRows = @(12,13,14,15,16,18,19)
foreach ($r in $runNoRows) {
    $ws.Range("D$r").Value = 4
}

# --- 3. Rows that stay "Failed" but get re-highlighted with a bold dark-red
#        font (new style). Build the format on B12 first, then fan it out to
#        the other affected rows via PasteSpecial so a single new style is
#        created and shared (matches the workbook's existing style-reuse
#        pattern).
$b12 = $ws.Range("B12")
$b12.Font.Bold = $true
$b12.Font.Color = 192
$b12.Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)

# --- 4. Clear the "Report Hide" manual-status note on row 17 (G17 becomes
#        an empty, but still styled/bordered, cell).
$ws.Range("G17").ClearContents()

# --- 5. Restore the saved cursor position to I19.
$null = $ws.Range("I19").Select()
